$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B width in the target file is 15.42578125 (same as column A).
# The COM layer quantizes ColumnWidth to 1/6-character steps when writing the
# stored <col width> attribute, so we pick the input that lands on the closest
# achievable stored width to the target.
$ws.Columns.Item(2).ColumnWidth = 14.667

$ws.Range("A1").Value = -0.083041962128703517
$ws.Range("B1").Value = 0.082651996587586041
$ws.Range("A2").Value = 0.0026297150996832386
$ws.Range("B2").Value = -0.0036128206759826753
$ws.Range("A3").Value = 0.10654769597175573
$ws.Range("B3").Value = -0.10718764326211172
$ws.Range("A4").Value = -0.17680974819446149
$ws.Range("B4").Value = 0.17571787601872657
$ws.Range("A5").Value = -0.16971787613331379
$ws.Range("B5").Value = 0.16751478353234361
$ws.Range("A6").Value = -0.054641627501935019
$ws.Range("B6").Value = 0.054603844907236354
$ws.Range("A7").Value = -0.039209218521671474
$ws.Range("B7").Value = 0.039135828441263953
$ws.Range("A8").Value = -0.019135828585696402
$ws.Range("B8").Value = 0.019110400371564751
$ws.Range("A9").Value = -0.013110400492032603
$ws.Range("B9").Value = 0.013096730614172714
$ws.Range("A10").Value = -0.0070967307350073838
$ws.Range("B10").Value = 0.0070967734081364142
$ws.Range("A11").Value = -0.0025967735263918712
$ws.Range("B11").Value = 0.0025946082326555597
$ws.Range("A12").Value = 0.003405391646460032
$ws.Range("B12").Value = -0.0034249701327571103
$ws.Range("A13").Value = 0.0094249700119828361
$ws.Range("B13").Value = -0.009438730007846452
$ws.Range("A14").Value = 0.021438729876695817
$ws.Range("B14").Value = -0.021481261365850379
$ws.Range("A15").Value = 0.02748126124556638
$ws.Range("B15").Value = -0.027551903592836879
$ws.Range("A16").Value = -0.015026548363602732
$ws.Range("B16").Value = 0.015004502885465953
$ws.Range("A17").Value = -0.0090045030057250841
$ws.Range("B17").Value = 0.0089999998743381937
$ws.Range("A18").Value = -0.075281711807804186
$ws.Range("B18").Value = 0.075183828927375629
$ws.Range("A19").Value = -0.027097273231270869
$ws.Range("B19").Value = 0.027014185748091002
$ws.Range("A20").Value = -0.018014185868283406
$ws.Range("B20").Value = 0.018004334673660694
$ws.Range("A21").Value = -0.0090043347940307328
$ws.Range("B21").Value = 0.0089999998794896285
$ws.Range("A22").Value = -0.093934875414641184
$ws.Range("B22").Value = 0.093625626978404597
$ws.Range("A23").Value = -0.084625627098326994
$ws.Range("B23").Value = 0.084125177295685027
$ws.Range("A24").Value = -0.042125177475573317
$ws.Range("B24").Value = 0.041999999819116951
$ws.Range("A25").Value = -0.094891580872868531
$ws.Range("B25").Value = 0.094648160292674532
$ws.Range("A26").Value = -0.088648160411999299
$ws.Range("B26").Value = 0.088335881190676702
$ws.Range("A27").Value = -0.082335881310648062
$ws.Range("B27").Value = 0.081273059129308933
$ws.Range("A28").Value = -0.07527305925199812
$ws.Range("B28").Value = 0.074536479481936091
$ws.Range("A29").Value = -0.062536479617065766
$ws.Range("B29").Value = 0.062170553279639051
$ws.Range("A30").Value = -0.042170553429759838
$ws.Range("B30").Value = 0.042020788387272745
$ws.Range("A31").Value = -0.027020788529965145
$ws.Range("B31").Value = 0.027001088011083141
$ws.Range("A32").Value = -0.0060010881645444769
$ws.Range("B32").Value = 0.0059999998727615633
